$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'23.873.56"
$ws.Range("E2").Value = "  -0.31%  "

# Row 3
$ws.Range("D3").Value = "'1.647.57"
$ws.Range("E3").Value = "  -0.06%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.91%  "

# Row 5
$ws.Range("D5").Value = "'310.42"
$ws.Range("E5").Value = "  +0.25%  "

# Row 6
$ws.Range("E6").Value = "  -0.75%  "

# Row 7
$ws.Range("D7").Value = "'0.3892"
$ws.Range("E7").Value = "  -0.75%  "

# Row 8
$ws.Range("D8").Value = "'0.3836"
$ws.Range("E8").Value = "  -1.11%  "

# Row 9
$ws.Range("D9").Value = "'51.07"
$ws.Range("E9").Value = "  -0.21%  "

# Row 10
$ws.Range("D10").Value = "'1.341"
$ws.Range("E10").Value = "  -1.43%  "

# Row 11
$ws.Range("D11").Value = "'1.001"
$ws.Range("E11").Value = "  -0.93%  "

# Row 12
$ws.Range("D12").Value = "'0.08432"
$ws.Range("E12").Value = "  -0.43%  "

# Row 13
$ws.Range("D13").Value = "'23.79"
$ws.Range("E13").Value = "  -0.41%  "

# Row 14
$ws.Range("D14").Value = "'7.014"
$ws.Range("E14").Value = "  -2.59%  "

# Row 15
$ws.Range("D15").Value = "'7.922"
$ws.Range("E15").Value = "  +0.50%  "

# Row 16
$ws.Range("D16").Value = "'0.00001312"
$ws.Range("E16").Value = "  +0.08%  "

# Row 17
$ws.Range("D17").Value = "'1.649.30"
$ws.Range("E17").Value = "  -0.10%  "

# Row 18
$ws.Range("D18").Value = "'93.90"
$ws.Range("E18").Value = "  -0.87%  "

# Row 19
$ws.Range("D19").Value = "'0.06958"
$ws.Range("E19").Value = "  -0.40%  "

# Row 20
$ws.Range("D20").Value = "'19.49"
$ws.Range("E20").Value = "  -2.80%  "

# Row 21
$ws.Range("D21").Value = "'6.935"
$ws.Range("E21").Value = "  +0.40%  "

# Row 22
$ws.Range("E22").Value = "  -0.71%  "

# Row 23
$ws.Range("E23").Value = "  -0.43%  "

# Row 24
$ws.Range("D24").Value = "'23.867.17"
$ws.Range("E24").Value = "  -0.35%  "

# Row 25
$ws.Range("D25").Value = "'2.440"
$ws.Range("E25").Value = "  -3.08%  "

# Row 26
$ws.Range("D26").Value = "'2.899"
$ws.Range("E26").Value = "  -4.73%  "

# Row 27
$ws.Range("D27").Value = "'21.92"
$ws.Range("E27").Value = "  -1.37%  "

# Row 28
$ws.Range("D28").Value = "'154.22"
$ws.Range("E28").Value = "  -0.73%  "

# Row 29
$ws.Range("D29").Value = "'5.378"
$ws.Range("E29").Value = "  +1.46%  "

# Row 30
$ws.Range("D30").Value = "'137.12"
$ws.Range("E30").Value = "  -2.01%  "

# Row 31
$ws.Range("D31").Value = "'7.691"
$ws.Range("E31").Value = "  -1.84%  "

# Row 32
$ws.Range("D32").Value = "'2.485"
$ws.Range("E32").Value = "  -1.45%  "

# Row 33
$ws.Range("D33").Value = "'1.831.01"
$ws.Range("E33").Value = "  -0.29%  "

# Row 34
$ws.Range("D34").Value = "'0.08108"
$ws.Range("E34").Value = "  +0.07%  "

# Row 35
$ws.Range("D35").Value = "'0.9870"
$ws.Range("E35").Value = "  -4.13%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02915"
$ws.Range("E36").Value = "  -3.40%  "

# Row 37
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'6.675"
$ws.Range("E37").Value = "  -0.52%  "

# Row 38
$ws.Range("D38").Value = "'0.2677"
$ws.Range("E38").Value = "  -1.21%  "

# Row 39
$ws.Range("D39").Value = "'10.47"
$ws.Range("E39").Value = "  -3.65%  "

# Row 40
$ws.Range("D40").Value = "'0.09110"
$ws.Range("E40").Value = "  -0.70%  "

# Row 41
$ws.Range("D41").Value = "'0.7537"
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("D42").Value = "'13.39"
$ws.Range("E42").Value = "  -0.92%  "

# Row 43
$ws.Range("E43").Value = "  -0.09%  "

# Row 44
$ws.Range("D44").Value = "'16.63"
$ws.Range("E44").Value = "  +2.88%  "

# Row 45
$ws.Range("D45").Value = "'0.6911"
$ws.Range("E45").Value = "  -0.24%  "

# Row 46
$ws.Range("D46").Value = "'2.434"
$ws.Range("E46").Value = "  -1.75%  "

# Row 47
$ws.Range("D47").Value = "'4.093"
$ws.Range("E47").Value = "  +0.12%  "

# Row 48
$ws.Range("E48").Value = "  -0.47%  "

# Row 49
$ws.Range("D49").Value = "'0.08263"
$ws.Range("E49").Value = "  -0.26%  "

# Row 50
$ws.Range("D50").Value = "'134.31"
$ws.Range("E50").Value = "  +0.15%  "

# Row 51
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'1.265"
$ws.Range("E51").Value = "  +6.20%  "
